# Apply the edits described by the commit:
#  - Replace three "OK..." model codes with "TEST..." placeholders
#    (C2: OK1013 -> TEST1013, C6: OK9100B -> TEST9100, C7: OK1224 -> TEST1224)
#  - Turn wrap-text on for D7 (the model-description cell for the new TEST1224 row)
#  - Move the sheet's active selection from C6 to H2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three product-code cells ---
$ws.Range("C2").Value = "TEST1013"
$ws.Range("C6").Value = "TEST9100"
$ws.Range("C7").Value = "TEST1224"

# --- Enable wrap text on D7 (model description for row 7) ---
$ws.Range("D7").WrapText = $true

# --- Update the selected cell shown when the sheet is reopened ---
$ws.Range("H2").Select()
